$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 42611.885740740741
$ws.Range("B8").Value = -2
$ws.Range("C8").Value = 54
$ws.Range("D8").Value = 44
$ws.Range("E8").Value = 33
$ws.Range("F8").Value = 66
$ws.Range("G8").Value = 17470
$ws.Range("H8").Value = 14211
$ws.Range("I8").Value = 905
$ws.Range("J8").Value = 135
$ws.Range("K8").Value = 112
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 8
$ws.Range("N8").Value = "Named"
